$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-05-09 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-10 Saturday", 2) | Out-Null
$d.Content.Find.Execute("43×90=3870", $true, $false, $false, $false, $false, $true, 1, $false, "22×93=2046", 2) | Out-Null
$d.Content.Find.Execute("43×14=602", $true, $false, $false, $false, $false, $true, 1, $false, "69×76=5244", 2) | Out-Null
$d.Content.Find.Execute("80×56=4480", $true, $false, $false, $false, $false, $true, 1, $false, "65×64=4160", 2) | Out-Null
$d.Content.Find.Execute("68×65=4420", $true, $false, $false, $false, $false, $true, 1, $false, "58×29=1682", 2) | Out-Null
$d.Content.Find.Execute("99×21=2079", $true, $false, $false, $false, $false, $true, 1, $false, "74×55=4070", 2) | Out-Null
$d.Content.Find.Execute("27×51=1377", $true, $false, $false, $false, $false, $true, 1, $false, "66×89=5874", 2) | Out-Null
$d.Content.Find.Execute("41×64=2624", $true, $false, $false, $false, $false, $true, 1, $false, "37×40=1480", 2) | Out-Null
$d.Content.Find.Execute("72×57=4104", $true, $false, $false, $false, $false, $true, 1, $false, "90×72=6480", 2) | Out-Null
$d.Content.Find.Execute("74×16=1184", $true, $false, $false, $false, $false, $true, 1, $false, "27×58=1566", 2) | Out-Null
$d.Content.Find.Execute("53×39=2067", $true, $false, $false, $false, $false, $true, 1, $false, "78×78=6084", 2) | Out-Null
$d.Content.Find.Execute("65×93=6045", $true, $false, $false, $false, $false, $true, 1, $false, "61×97=5917", 2) | Out-Null
$d.Content.Find.Execute("22×70=1540", $true, $false, $false, $false, $false, $true, 1, $false, "88×16=1408", 2) | Out-Null
$d.Content.Find.Execute("33×46=1518", $true, $false, $false, $false, $false, $true, 1, $false, "76×11=836", 2) | Out-Null
$d.Content.Find.Execute("18×46=828", $true, $false, $false, $false, $false, $true, 1, $false, "72×35=2520", 2) | Out-Null
$d.Content.Find.Execute("43×82=3526", $true, $false, $false, $false, $false, $true, 1, $false, "88×40=3520", 2) | Out-Null
$d.Content.Find.Execute("14×17=238", $true, $false, $false, $false, $false, $true, 1, $false, "32×30=960", 2) | Out-Null
$d.Content.Find.Execute("55×20=1100", $true, $false, $false, $false, $false, $true, 1, $false, "94×84=7896", 2) | Out-Null
$d.Content.Find.Execute("79×79=6241", $true, $false, $false, $false, $false, $true, 1, $false, "22×51=1122", 2) | Out-Null
$d.Content.Find.Execute("27×11=297", $true, $false, $false, $false, $false, $true, 1, $false, "62×37=2294", 2) | Out-Null
$d.Content.Find.Execute("61×74=4514", $true, $false, $false, $false, $false, $true, 1, $false, "24×81=1944", 2) | Out-Null
$d.Content.Find.Execute("44×53=2332", $true, $false, $false, $false, $false, $true, 1, $false, "72×18=1296", 2) | Out-Null
$d.Content.Find.Execute("48×68=3264", $true, $false, $false, $false, $false, $true, 1, $false, "12×17=204", 2) | Out-Null
$d.Content.Find.Execute("59×73=4307", $true, $false, $false, $false, $false, $true, 1, $false, "59×79=4661", 2) | Out-Null
$d.Content.Find.Execute("59×75=4425", $true, $false, $false, $false, $false, $true, 1, $false, "75×47=3525", 2) | Out-Null
$d.Content.Find.Execute("89×68=6052", $true, $false, $false, $false, $false, $true, 1, $false, "50×16=800", 2) | Out-Null
